# Update countries & provincias Spain
# Applies the 2020-04-03 20:50 -> 21:25 data refresh to the "Pais" sheet:
#  - a handful of country rows get updated case/death/recovered counts
#  - some countries leapfrog each other in the (rank-sorted) table, so the
#    country name in column A for a given row changes along with its data
#  - the "last updated" footer timestamp text changes

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Row 4: Estados Unidos
$ws.Cells.Item(4, 1).Value = "Estados Unidos"
$ws.Cells.Item(4, 2).Value = 266558
$ws.Cells.Item(4, 3).Value = 21681
$ws.Cells.Item(4, 4).Value = 11983
$ws.Cells.Item(4, 5).Value = 247772
$ws.Cells.Item(4, 6).Value = 5781
$ws.Cells.Item(4, 7).Value = 733
$ws.Cells.Item(4, 8).Value = 6803

# Row 7: Alemania
$ws.Cells.Item(7, 1).Value = "Alemania"
$ws.Cells.Item(7, 2).Value = 90964
$ws.Cells.Item(7, 3).Value = 6170
$ws.Cells.Item(7, 4).Value = 24575
$ws.Cells.Item(7, 5).Value = 65155
$ws.Cells.Item(7, 6).Value = 3936
$ws.Cells.Item(7, 7).Value = 127
$ws.Cells.Item(7, 8).Value = 1234

# Row 16: Canada
$ws.Cells.Item(16, 1).Value = "Canada"
$ws.Cells.Item(16, 2).Value = 12375
$ws.Cells.Item(16, 3).Value = 1092
$ws.Cells.Item(16, 4).Value = 1979
$ws.Cells.Item(16, 5).Value = 10218
$ws.Cells.Item(16, 6).Value = 120
$ws.Cells.Item(16, 7).Value = 5
$ws.Cells.Item(16, 8).Value = 178

# Row 20: Brasil
$ws.Cells.Item(20, 1).Value = "Brasil"
$ws.Cells.Item(20, 2).Value = 8261
$ws.Cells.Item(20, 3).Value = 217
$ws.Cells.Item(20, 4).Value = 127
$ws.Cells.Item(20, 5).Value = 7790
$ws.Cells.Item(20, 6).Value = 296
$ws.Cells.Item(20, 7).Value = 20
$ws.Cells.Item(20, 8).Value = 344

# Row 23: Noruega
$ws.Cells.Item(23, 1).Value = "Noruega"
$ws.Cells.Item(23, 2).Value = 5370
$ws.Cells.Item(23, 3).Value = 223
$ws.Cells.Item(23, 4).Value = 32
$ws.Cells.Item(23, 5).Value = 5279
$ws.Cells.Item(23, 6).Value = 96
$ws.Cells.Item(23, 7).Value = 9
$ws.Cells.Item(23, 8).Value = 59

# Row 24: Australia
$ws.Cells.Item(24, 1).Value = "Australia"
$ws.Cells.Item(24, 2).Value = 5350
$ws.Cells.Item(24, 3).Value = 36
$ws.Cells.Item(24, 4).Value = 585
$ws.Cells.Item(24, 5).Value = 4737
$ws.Cells.Item(24, 6).Value = 85
$ws.Cells.Item(24, 7).Value = 3
$ws.Cells.Item(24, 8).Value = 28

# Row 58: Ucrania
$ws.Cells.Item(58, 1).Value = "Ucrania"
$ws.Cells.Item(58, 2).Value = 1072
$ws.Cells.Item(58, 3).Value = 175
$ws.Cells.Item(58, 4).Value = 22
$ws.Cells.Item(58, 5).Value = 1027
$ws.Cells.Item(58, 6).Value = 16
$ws.Cells.Item(58, 7).Value = 1
$ws.Cells.Item(58, 8).Value = 23

# Row 59: Egipto
$ws.Cells.Item(59, 1).Value = "Egipto"
$ws.Cells.Item(59, 2).Value = 985
$ws.Cells.Item(59, 3).Value = 120
$ws.Cells.Item(59, 4).Value = 216
$ws.Cells.Item(59, 5).Value = 703
$ws.Cells.Item(59, 6).Value = 0
$ws.Cells.Item(59, 7).Value = 8
$ws.Cells.Item(59, 8).Value = 66

# Row 60: Estonia
$ws.Cells.Item(60, 1).Value = "Estonia"
$ws.Cells.Item(60, 2).Value = 961
$ws.Cells.Item(60, 3).Value = 103
$ws.Cells.Item(60, 4).Value = 48
$ws.Cells.Item(60, 5).Value = 901
$ws.Cells.Item(60, 6).Value = 16
$ws.Cells.Item(60, 7).Value = 1
$ws.Cells.Item(60, 8).Value = 12

# Row 118: Mayotte
$ws.Cells.Item(118, 1).Value = "Mayotte"
$ws.Cells.Item(118, 2).Value = 128
$ws.Cells.Item(118, 3).Value = 12
$ws.Cells.Item(118, 4).Value = 10
$ws.Cells.Item(118, 5).Value = 116
$ws.Cells.Item(118, 6).Value = 3
$ws.Cells.Item(118, 7).Value = 1
$ws.Cells.Item(118, 8).Value = 2

# Row 119: Kenia
$ws.Cells.Item(119, 1).Value = "Kenia"
$ws.Cells.Item(119, 2).Value = 122
$ws.Cells.Item(119, 3).Value = 12
$ws.Cells.Item(119, 4).Value = 4
$ws.Cells.Item(119, 5).Value = 114
$ws.Cells.Item(119, 6).Value = 2
$ws.Cells.Item(119, 7).Value = 1
$ws.Cells.Item(119, 8).Value = 4

# Row 137: Uganda
$ws.Cells.Item(137, 1).Value = "Uganda"
$ws.Cells.Item(137, 2).Value = 48
$ws.Cells.Item(137, 3).Value = 3
$ws.Cells.Item(137, 4).Value = 0
$ws.Cells.Item(137, 5).Value = 48
$ws.Cells.Item(137, 6).Value = 0
$ws.Cells.Item(137, 7).Value = 0
$ws.Cells.Item(137, 8).Value = 0

# Row 138: Jamaica
$ws.Cells.Item(138, 1).Value = "Jamaica"
$ws.Cells.Item(138, 2).Value = 47
$ws.Cells.Item(138, 3).Value = 0
$ws.Cells.Item(138, 4).Value = 2
$ws.Cells.Item(138, 5).Value = 42
$ws.Cells.Item(138, 6).Value = 0
$ws.Cells.Item(138, 7).Value = 0
$ws.Cells.Item(138, 8).Value = 3

# Row 139: El Salvador
$ws.Cells.Item(139, 1).Value = "El Salvador"
$ws.Cells.Item(139, 2).Value = 46
$ws.Cells.Item(139, 3).Value = 5
$ws.Cells.Item(139, 4).Value = 0
$ws.Cells.Item(139, 5).Value = 44
$ws.Cells.Item(139, 6).Value = 4
$ws.Cells.Item(139, 7).Value = 0
$ws.Cells.Item(139, 8).Value = 2

# Footer timestamp text (row 1)
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 3 de Abril de 2020 a las 21:25"
